# Weekly update: a new daily price record (2021-10-13) was added to the
# "Naranja" dataset. It is inserted as a new row 30, pushing the existing
# rows 30-68 down to 31-69 (dimension grows from A1:T68 to A1:T69).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 30, shifting rows 30:68 -> 31:69.
$ws.Rows("30:30").Insert()

# Populate the newly inserted row 30 with the new record.
$ws.Range("A30").Value = 1
$ws.Range("B30").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C30").Value = "Arica y Parinacota"
$ws.Range("D30").Value = 44482
$ws.Range("E30").Value = 15
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100102
$ws.Range("H30").Value = "Cítricos"
$ws.Range("I30").Value = 100102005
$ws.Range("J30").Value = "Naranja"
$ws.Range("K30").Value = "Lane Late"
$ws.Range("L30").Value = "Segunda"
$ws.Range("M30").Value = 300
$ws.Range("N30").Value = 650
$ws.Range("O30").Value = 700
$ws.Range("P30").Value = 675
$ws.Range("Q30").Value = "`$/kilo (en caja de 20 kilos)"
$ws.Range("R30").Value = "Región de Coquimbo"
$ws.Range("S30").Value = 675
$ws.Range("T30").Value = 1
